{"js": "// Target change (see XML diff in the task):\n//   - word/document.xml: the xmlns:* declarations on the root <w:document>\n//     element are re-ordered alphabetically (mc:Ignorable moves after them),\n//     and the attributes of <w:pgSz>/<w:pgMar> inside <w:sectPr> are\n//     re-ordered alphabetically.\n//   - word/styles.xml: attributes on <w:rFonts>, <w:lang>, <w:latentStyles>,\n//     every <w:lsdException>, every <w:style> (incl. w:tblInd / w:tblCellMar\n//     children of the \"TableauNormal\" style), are likewise re-ordered\n//     alphabetically.\n//\n// Every one of those hunks is a pure XML-attribute-order change: the\n// element names, attribute names, attribute values, text runs, formatting\n// and document structure are byte-for-byte identical before and after\n// (this is confirmed by sorting the attributes of the \"before\" OOXML and\n// comparing it to the \"after\" OOXML - they match exactly). The commit\n// message confirms this is a tooling change (\"Fixed POI packaging and\n// upgraded to POI 3.15\"): re-serializing the parts with a newer Apache POI\n// version happens to emit attributes in a different (alphabetical) order.\n//\n// Attribute serialization order is not part of the Word document object\n// model exposed by Office.js (Word.* API): there is no property anywhere\n// on Document/Body/Section/PageSetup/Style/Font that lets a script choose\n// the order in which the underlying XML writer emits an element's\n// attributes, and re-setting those properties to their current values does\n// not change how the host re-serializes them. So there is no Office.js\n// call that reproduces this particular diff, and no visible/semantic\n// document content needs to change to match it (the rendered document,\n// its text, formatting, styles and page setup values are unaffected).\n//\n// Intentionally a no-op edit: we only load/sync so the script is a valid,\n// side-effect-free pass-through rather than risk unrelated changes (e.g.\n// touching styles/sections has been observed to mint unrelated namespace\n// declarations / rewrite unrelated runs as a side effect in this host).\ncontext.document.body.load(\"text\");\nawait context.sync();\n", "ps1": "# Target change (see XML diff in the task):\n#   - word/document.xml: the xmlns:* declarations on the root <w:document>\n#     element are re-ordered alphabetically (mc:Ignorable moves after them),\n#     and the attributes of <w:pgSz>/<w:pgMar> inside <w:sectPr> are\n#     re-ordered alphabetically.\n#   - word/styles.xml: attributes on <w:rFonts>, <w:lang>, <w:latentStyles>,\n#     every <w:lsdException>, every <w:style> (incl. w:tblInd / w:tblCellMar\n#     children of the \"TableauNormal\" style), are likewise re-ordered\n#     alphabetically.\n#\n# Every one of those hunks is a pure XML-attribute-order change: the\n# element names, attribute names, attribute values, text runs, formatting\n# and document structure are byte-for-byte identical before and after\n# (sorting the attributes of the \"before\" OOXML alphabetically reproduces\n# the \"after\" OOXML exactly). The commit message confirms this is a\n# tooling change (\"Fixed POI packaging and upgraded to POI 3.15\"): saving\n# the parts with a newer Apache POI version happens to emit attributes in\n# a different (alphabetical) order.\n#\n# Attribute serialization order is not part of the Word COM object model:\n# there is no property anywhere on Document/Sections/PageSetup/Styles/Font\n# that lets a script choose the order in which the underlying XML writer\n# emits an element's attributes, and re-assigning those properties to\n# their current values (PageSetup.TopMargin/PageWidth/..., Styles(...).*,\n# etc.) does not change how the host re-serializes them - it was verified\n# to keep writing attributes in the same fixed order regardless. So there\n# is no COM call sequence that reproduces this particular diff, and no\n# visible/semantic document content needs to change to match it (the\n# rendered document, its text, formatting, styles and page setup values\n# are unaffected).\n#\n# Intentionally a no-op edit: we only touch $d.Content read-only so the\n# script is a valid, side-effect-free pass-through rather than risk\n# unrelated changes (touching PageSetup/Styles properties has been\n# observed to mint unrelated namespace declarations / rewrite unrelated\n# field-code runs as a side effect in this host).\n$d = $word.ActiveDocument\n$null = $d.Content\n"}
